$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.066.10"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "2.349.89"
$ws.Range("E3").Value = "  -5.41%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'474.28"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "'144.46"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.611"
$ws.Range("E8").Value = "  +19.52%  "
$ws.Range("E9").Value = "  -5.28%  "
$ws.Range("D10").Value = "'0.0970"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  -6.00%  "
$ws.Range("D12").Value = "'0.323"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D15").Value = "55.035.72"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").Value = "'20.00"
$ws.Range("E16").Value = "  -5.12%  "
$ws.Range("D17").Value = "'0.0000130"
$ws.Range("E17").Value = "  -4.08%  "
$ws.Range("D18").Value = "2.356.89"
$ws.Range("E18").Value = "  -5.21%  "
$ws.Range("D19").Value = "'4.58"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "'315.99"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'9.59"
$ws.Range("E21").Value = "  -4.91%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D23").Value = "'5.62"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").Value = "'56.89"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'0.393"
$ws.Range("E26").Value = "  -4.41%  "
$ws.Range("E27").Value = "  -5.36%  "
$ws.Range("D28").Value = "2.451.58"
$ws.Range("E28").Value = "  -5.18%  "
$ws.Range("D29").Value = "'7.14"
$ws.Range("E29").Value = "  -6.58%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "0.0₃0751"
$ws.Range("E31").Value = "  -5.59%  "
$ws.Range("D32").Value = "'146.61"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").Value = "'18.20"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").Value = "'5.07"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("D36").Value = "'3.58"
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("D38").Value = "'0.811"
$ws.Range("E38").Value = "  -6.14%  "
$ws.Range("D39").Value = "'33.70"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").Value = "'0.0992"
$ws.Range("E40").Value = "  +7.07%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'3.40"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").Value = "'0.0518"
$ws.Range("E45").Value = "  -6.97%  "
$ws.Range("D46").Value = "'10.17"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'251.46"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "'0.0220"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("D49").Value = "'4.35"
$ws.Range("E49").Value = "  -9.01%  "
$ws.Range("D50").Value = "'16.72"
$ws.Range("E50").Value = "  -4.96%  "
$ws.Range("D51").Value = "1.777.29"
$ws.Range("E51").Value = "  -5.03%  "
